$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.164110538044052
$ws.Range("C2").Value = 0.3990739875612785
$ws.Range("D2").Value = 0.2717830387705941
$ws.Range("E2").Value = 0.5213281488377489
$ws.Range("F2").Value = 0.5059450074181695
$ws.Range("G2").Value = 23

# Row 3
$ws.Range("B3").Value = 0.05421946485828021
$ws.Range("C3").Value = 0.3743046632812781
$ws.Range("D3").Value = 0.2039986448588482
$ws.Range("E3").Value = 0.4516620914564872
$ws.Range("F3").Value = 0.4589478484100895
$ws.Range("G3").Value = 22

# Row 4
$ws.Range("B4").Value = 0.2108166877076684
$ws.Range("C4").Value = 0.3911061312796325
$ws.Range("D4").Value = 0.3070625528270741
$ws.Range("E4").Value = 0.5541322521087129
$ws.Range("F4").Value = 0.5251188635552846
$ws.Range("G4").Value = 21

# Row 5
$ws.Range("B5").Value = 0.1396966153068137
$ws.Range("C5").Value = 0.3936428606436134
$ws.Range("D5").Value = 0.2364141445760344
$ws.Range("E5").Value = 0.4862243767809615
$ws.Range("F5").Value = 0.477822913957676
$ws.Range("G5").Value = 20

# Row 6
$ws.Range("B6").Value = 0.2055713184881964
$ws.Range("C6").Value = 0.3985724376705788
$ws.Range("D6").Value = 0.2720016382854029
$ws.Range("E6").Value = 0.521537763048279
$ws.Range("F6").Value = 0.4924484944702406
$ws.Range("G6").Value = 19

# Row 7
$ws.Range("B7").Value = 0.1942106637073617
$ws.Range("C7").Value = 0.387649315436281
$ws.Range("D7").Value = 0.2468445352345524
$ws.Range("E7").Value = 0.4968345149388803
$ws.Range("F7").Value = 0.4705617144037521
$ws.Range("G7").Value = 18

# Row 8
$ws.Range("B8").Value = 0.2291527386783319
$ws.Range("C8").Value = 0.4681846348907073
$ws.Range("D8").Value = 0.4295473340169714
$ws.Range("E8").Value = 0.6553986069690501
$ws.Range("F8").Value = 0.6329305875422011
$ws.Range("G8").Value = 17

$wb.Save()
